$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# The orchestrator now writes a per-line Status/Notes back onto the order
# sheet once each product line has been processed against the catalog and,
# on success, the confirmation e-mail has gone out.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = "OK"
    $ws.Cells.Item($r, 4).Value = "Order Number: 512464"
}

# Row 9 references "Ipoh Coff", which isn't a catalog product (the real
# item is "Ipoh Coffee"), so the orchestrator reports it as an error
# instead of a successful order confirmation.
$ws.Cells.Item(9, 3).Value = "Erro"
$ws.Cells.Item(9, 4).Value = "Product Ipoh Coff not Found"

# Column E was just a blank placeholder left over from the old config-file
# workflow; it's no longer produced by the orchestrator, so drop it.
$ws.Columns.Item(5).Delete()

$ws.Range("C9").Select()
